# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit recalculation updates across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 294.14285
$ws.Range("I5").Value = 221.4
$ws.Range("K5").Value = 221.4
$ws.Range("M5").Value = -106.4
$ws.Range("H28").Value = 556569.1
$ws.Range("I28").Value = 741374.1
$ws.Range("J28").Value = 2154
$ws.Range("K28").Value = 741374.1
$ws.Range("L28").Value = 2154
$ws.Range("M28").Value = -740889.1
$ws.Range("N28").Value = -3124
$ws.Range("H33").Value = 154.56
$ws.Range("I33").Value = 160.66667
$ws.Range("K33").Value = 160.66667
$ws.Range("M33").Value = 68.33332999999999
$ws.Range("H64").Value = 5267.7646
$ws.Range("I64").Value = 4137
$ws.Range("J64").Value = 6883.143
$ws.Range("K64").Value = 4137
$ws.Range("L64").Value = 6883.143
$ws.Range("M64").Value = -3889
$ws.Range("N64").Value = -7379.143
$ws.Range("H67").Value = 5267.7646
$ws.Range("I67").Value = 4137
$ws.Range("J67").Value = 6883.143
$ws.Range("K67").Value = 4137
$ws.Range("L67").Value = 6883.143
$ws.Range("M67").Value = -3279
$ws.Range("N67").Value = -8599.143
$ws.Range("H75").Value = 106000
$ws.Range("J75").Value = 144000
$ws.Range("L75").Value = 144000
$ws.Range("N75").Value = -145872
$ws.Range("H76").Value = 3475315.5
$ws.Range("I76").Value = 4447284
$ws.Range("K76").Value = 4447284
$ws.Range("M76").Value = -4446969
$ws.Range("H78").Value = 106000
$ws.Range("J78").Value = 144000
$ws.Range("L78").Value = 432000
$ws.Range("N78").Value = -441360
$ws.Range("H79").Value = 3475315.5
$ws.Range("I79").Value = 4447284
$ws.Range("K79").Value = 4447284
$ws.Range("M79").Value = -4446192
$ws.Range("H93").Value = 25933.572
$ws.Range("J93").Value = 25933.572
$ws.Range("L93").Value = 25933.572
$ws.Range("N93").Value = -30925.572
$ws.Range("H96").Value = 638.7857
$ws.Range("I96").Value = 649.46155
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 1948.38465
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -575.38465
$ws.Range("N96").Value = -4246
$ws.Range("H98").Value = 363930.2
$ws.Range("I98").Value = 388749.53
$ws.Range("J98").Value = 4050
$ws.Range("K98").Value = 388749.53
$ws.Range("L98").Value = 4050
$ws.Range("M98").Value = -387251.53
$ws.Range("N98").Value = -7046
$ws.Range("H107").Value = 529532.9399999999
$ws.Range("I107").Value = 585183.2
$ws.Range("J107").Value = 856
$ws.Range("K107").Value = 585183.2
$ws.Range("L107").Value = 856
$ws.Range("M107").Value = -583263.2
$ws.Range("N107").Value = -4696
$ws.Range("H112").Value = 8524013
$ws.Range("I112").Value = 589.5
$ws.Range("J112").Value = 9092241
$ws.Range("K112").Value = 1768.5
$ws.Range("L112").Value = 27276723
$ws.Range("M112").Value = -660.5
$ws.Range("N112").Value = -27278939
$ws.Range("H122").Value = 363930.2
$ws.Range("I122").Value = 388749.53
$ws.Range("J122").Value = 4050
$ws.Range("K122").Value = 1166248.59
$ws.Range("L122").Value = 12150
$ws.Range("M122").Value = -1163798.59
$ws.Range("N122").Value = -17050
$ws.Range("H126").Value = 34945
$ws.Range("J126").Value = 34945
$ws.Range("L126").Value = 34945
$ws.Range("N126").Value = -44825
$ws.Range("H127").Value = 1059.1875
$ws.Range("I127").Value = 626.1667
$ws.Range("J127").Value = 1319
$ws.Range("K127").Value = 1878.5001
$ws.Range("L127").Value = 3957
$ws.Range("M127").Value = 3081.4999
$ws.Range("N127").Value = -13877
$ws.Range("H129").Value = 1279.3043
$ws.Range("J129").Value = 1487.1052
$ws.Range("L129").Value = 4461.3156
$ws.Range("N129").Value = -14461.3156
$ws.Range("H135").Value = 936.8182
$ws.Range("I135").Value = 936.8182
$ws.Range("K135").Value = 8431.363800000001
$ws.Range("M135").Value = -5896.363800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 59750
$ws.Range("J24").Value = 59750
$ws.Range("L24").Value = 59750
$ws.Range("N24").Value = -60498
$ws.Range("H61").Value = 2357.96
$ws.Range("I61").Value = 1705.4
$ws.Range("J61").Value = 4968.2
$ws.Range("K61").Value = 1705.4
$ws.Range("L61").Value = 4968.2
$ws.Range("M61").Value = -1493.4
$ws.Range("N61").Value = -5392.2
$ws.Range("H100").Value = 59750
$ws.Range("J100").Value = 59750
$ws.Range("L100").Value = 59750
$ws.Range("N100").Value = -61914
$ws.Range("H102").Value = 1859.875
$ws.Range("I102").Value = 1579.8334
$ws.Range("K102").Value = 1579.8334
$ws.Range("M102").Value = 42.16660000000002
$ws.Range("H132").Value = 3151.2432
$ws.Range("I132").Value = 2577.3103
$ws.Range("J132").Value = 5231.75
$ws.Range("K132").Value = 7731.9309
$ws.Range("L132").Value = 15695.25
$ws.Range("M132").Value = -5201.9309
$ws.Range("N132").Value = -20755.25
$ws.Range("H136").Value = 2357.96
$ws.Range("I136").Value = 1705.4
$ws.Range("J136").Value = 4968.2
$ws.Range("K136").Value = 5116.200000000001
$ws.Range("L136").Value = 14904.6
$ws.Range("M136").Value = -2566.200000000001
$ws.Range("N136").Value = -20004.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 50000
$ws.Range("I19").Value = 50000
$ws.Range("K19").Value = 50000
$ws.Range("M19").Value = -49827
$ws.Range("H99").Value = 1778.091
$ws.Range("I99").Value = 1444.875
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 1444.875
$ws.Range("L99").Value = 2666.6667
$ws.Range("M99").Value = 53.125
$ws.Range("N99").Value = -5662.6667
$ws.Range("H105").Value = 2948.64
$ws.Range("I105").Value = 2740.2632
$ws.Range("J105").Value = 3608.5
$ws.Range("K105").Value = 2740.2632
$ws.Range("L105").Value = 3608.5
$ws.Range("M105").Value = -993.2631999999999
$ws.Range("N105").Value = -7102.5
$ws.Range("H134").Value = 3158.4722
$ws.Range("I134").Value = 2276.8147
$ws.Range("K134").Value = 6830.4441
$ws.Range("M134").Value = -4295.4441

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1580
$ws.Range("I16").Value = 1401.4286
$ws.Range("K16").Value = 1401.4286
$ws.Range("M16").Value = -1114.4286
$ws.Range("H33").Value = 37466.6
$ws.Range("I33").Value = 37466.6
$ws.Range("K33").Value = 37466.6
$ws.Range("M33").Value = -37087.6
$ws.Range("H52").Value = 44833.332
$ws.Range("J52").Value = 44833.332
$ws.Range("L52").Value = 44833.332
$ws.Range("N52").Value = -45421.332
$ws.Range("H62").Value = 20124.285
$ws.Range("I62").Value = 23567.273
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 23567.273
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -22943.273
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 20124.285
$ws.Range("I65").Value = 23567.273
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 117836.365
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -114716.365
$ws.Range("N65").Value = -43740
$ws.Range("H113").Value = 1580
$ws.Range("I113").Value = 1401.4286
$ws.Range("K113").Value = 1401.4286
$ws.Range("M113").Value = 768.5714
$ws.Range("H132").Value = 1791
$ws.Range("I132").Value = 1245.5
$ws.Range("K132").Value = 3736.5
$ws.Range("M132").Value = -1206.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 1032.5714
$ws.Range("I119").Value = 738
$ws.Range("J119").Value = 2800
$ws.Range("K119").Value = 2214
$ws.Range("L119").Value = 8400
$ws.Range("M119").Value = 2624
$ws.Range("N119").Value = -18076

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5644.727
$ws.Range("I70").Value = 5567
$ws.Range("J70").Value = 6080
$ws.Range("K70").Value = 5567
$ws.Range("L70").Value = 6080
$ws.Range("M70").Value = -5297
$ws.Range("N70").Value = -6620
$ws.Range("H73").Value = 5644.727
$ws.Range("I73").Value = 5567
$ws.Range("J73").Value = 6080
$ws.Range("K73").Value = 5567
$ws.Range("L73").Value = 6080
$ws.Range("M73").Value = -4631
$ws.Range("N73").Value = -7952
$ws.Range("H93").Value = 28938
$ws.Range("J93").Value = 28938
$ws.Range("L93").Value = 28938
$ws.Range("N93").Value = -32682
$ws.Range("H122").Value = 1761.75
$ws.Range("I122").Value = 907
$ws.Range("J122").Value = 1883.8572
$ws.Range("K122").Value = 2721
$ws.Range("L122").Value = 5651.571599999999
$ws.Range("M122").Value = -271
$ws.Range("N122").Value = -10551.5716
$ws.Range("H132").Value = 3205.55
$ws.Range("I132").Value = 2335.9666
$ws.Range("J132").Value = 5814.3
$ws.Range("K132").Value = 7007.899800000001
$ws.Range("L132").Value = 17442.9
$ws.Range("M132").Value = -4477.899800000001
$ws.Range("N132").Value = -22502.9
$ws.Range("H136").Value = 19311.6
$ws.Range("J136").Value = 19103.715
$ws.Range("L136").Value = 57311.145
$ws.Range("N136").Value = -62411.145
$ws.Range("H137").Value = 50169.5
$ws.Range("J137").Value = 50169.5
$ws.Range("L137").Value = 50169.5
$ws.Range("N137").Value = -60369.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 31750
$ws.Range("J96").Value = 31750
$ws.Range("L96").Value = 31750
$ws.Range("N96").Value = -37242

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 17766.666
$ws.Range("I9").Value = 25250
$ws.Range("J9").Value = 2800
$ws.Range("K9").Value = 25250
$ws.Range("L9").Value = 2800
$ws.Range("M9").Value = -25110
$ws.Range("N9").Value = -3080
$ws.Range("H32").Value = 33014.5
$ws.Range("I32").Value = 6000
$ws.Range("K32").Value = 6000
$ws.Range("M32").Value = -5683
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H100").Value = 866.5
$ws.Range("I100").Value = 1150
$ws.Range("J100").Value = 724.75
$ws.Range("K100").Value = 2300
$ws.Range("L100").Value = 1449.5
$ws.Range("M100").Value = -1759
$ws.Range("N100").Value = -2531.5
